$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new headers I0 (col I) and IF (col J) in the header row, matching the
# style already used by the other header cells (e.g. H1).
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy() | Out-Null
$ws.Range("I1:J1").PasteSpecial(-4122) | Out-Null

# Fill in the I0/IF data for rows 2-77.
$i0Values = @(8,9,8,8,6,7,8,7,8,7,7,9,9,9,8,10,7,7,8,7,6,8,8,7,7,7,9,7,9,8,7,7,7,4,7,7,7,6,8,8,8,8,6,9,6,7,10,8,5,8,8,4,8,7,6,8,5,6,7,6,7,8,5,7,7,6,7,7,7,6,5,1,6,4,6,3)
$ifValues = @(8,9,8,8,6,7,8,7,8,7,7,9,9,9,8,10,7,7,8,8,6,8,8,7,7,7,9,7,9,8,7,7,7,4,7,7,7,6,9,8,8,8,7,9,6,7,10,8,6,8,8,5,8,7,6,8,5,6,7,6,7,8,6,7,7,7,7,7,7,7,5,1,6,4,6,3)

for ($idx = 0; $idx -lt $i0Values.Length; $idx++) {
    $row = $idx + 2
    $ws.Cells.Item($row, 9).Value = $i0Values[$idx]
    $ws.Cells.Item($row, 10).Value = $ifValues[$idx]
}

$excel.CutCopyMode = 0
